$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4042.6226
$ws.Range("I64").Value = 3872.05
$ws.Range("K64").Value = 3872.05
$ws.Range("M64").Value = -3624.05

$ws.Range("H67").Value = 4042.6226
$ws.Range("I67").Value = 3872.05
$ws.Range("K67").Value = 3872.05
$ws.Range("M67").Value = -3014.05

$ws.Range("H70").Value = 4264.357
$ws.Range("I70").Value = 1300
$ws.Range("J70").Value = 5072.8184
$ws.Range("K70").Value = 3900
$ws.Range("L70").Value = 15218.4552
$ws.Range("M70").Value = -3630
$ws.Range("N70").Value = -15758.4552

$ws.Range("H73").Value = 4264.357
$ws.Range("I73").Value = 1300
$ws.Range("J73").Value = 5072.8184
$ws.Range("K73").Value = 3900
$ws.Range("L73").Value = 15218.4552
$ws.Range("M73").Value = -2964
$ws.Range("N73").Value = -17090.4552

$ws.Range("H129").Value = 836.84906
$ws.Range("J129").Value = 901.5
$ws.Range("L129").Value = 2704.5
$ws.Range("N129").Value = -12704.5

$ws.Range("H138").Value = 2884.4385
$ws.Range("J138").Value = 3641.0527
$ws.Range("L138").Value = 10923.1581
$ws.Range("N138").Value = -21203.1581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 9785
$ws.Range("J24").Value = 9785
$ws.Range("L24").Value = 9785

$ws.Range("H100").Value = 9785
$ws.Range("J100").Value = 9785
$ws.Range("L100").Value = 9785

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 28762
$ws.Range("J93").Value = 28762
$ws.Range("L93").Value = 28762

$ws.Range("H122").Value = 40000
$ws.Range("J122").Value = 40000
$ws.Range("L122").Value = 40000

$ws.Range("H134").Value = 2251.139
$ws.Range("I134").Value = 1108.52
$ws.Range("K134").Value = 3325.56
$ws.Range("M134").Value = -790.5599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1644.3091
$ws.Range("I31").Value = 1331.5116
$ws.Range("J31").Value = 2765.1667
$ws.Range("K31").Value = 1331.5116
$ws.Range("L31").Value = 2765.1667
$ws.Range("M31").Value = -1036.5116
$ws.Range("N31").Value = -3355.1667

$ws.Range("H34").Value = 1644.3091
$ws.Range("I34").Value = 1331.5116
$ws.Range("J34").Value = 2765.1667
$ws.Range("K34").Value = 1331.5116
$ws.Range("L34").Value = 2765.1667
$ws.Range("M34").Value = -1129.5116
$ws.Range("N34").Value = -3169.1667

$ws.Range("H58").Value = 1629.6227
$ws.Range("I58").Value = 1147.7727
$ws.Range("J58").Value = 1971.5807
$ws.Range("K58").Value = 1147.7727
$ws.Range("L58").Value = 1971.5807
$ws.Range("M58").Value = -944.7727
$ws.Range("N58").Value = -2377.5807

$ws.Range("H99").Value = 7988
$ws.Range("I99").Value = 1380.6
$ws.Range("K99").Value = 1380.6
$ws.Range("M99").Value = 117.4000000000001

$ws.Range("H126").Value = 7988
$ws.Range("I126").Value = 1380.6
$ws.Range("K126").Value = 4141.799999999999
$ws.Range("M126").Value = -1671.799999999999

$ws.Range("H134").Value = 3209.2942
$ws.Range("I134").Value = 3134.1538
$ws.Range("J134").Value = 3453.5
$ws.Range("K134").Value = 9402.4614
$ws.Range("L134").Value = 10360.5
$ws.Range("M134").Value = -6867.4614
$ws.Range("N134").Value = -15430.5

$ws.Range("H136").Value = 1629.6227
$ws.Range("I136").Value = 1147.7727
$ws.Range("J136").Value = 1971.5807
$ws.Range("K136").Value = 3443.3181
$ws.Range("L136").Value = 5914.742099999999
$ws.Range("M136").Value = -893.3181
$ws.Range("N136").Value = -11014.7421

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 129.9375
$ws.Range("I8").Value = 129.9375
$ws.Range("K8").Value = 389.8125
$ws.Range("M8").Value = -250.8125

$ws.Range("H86").Value = 902.4
$ws.Range("J86").Value = 902.4
$ws.Range("L86").Value = 2707.2

$ws.Range("H89").Value = 902.4
$ws.Range("J89").Value = 902.4
$ws.Range("L89").Value = 8121.599999999999

$ws.Range("H106").Value = 7929.5
$ws.Range("J106").Value = 7929.5
$ws.Range("L106").Value = 23788.5
$ws.Range("N106").Value = -25680.5

$ws.Range("H112").Value = 1440.4517
$ws.Range("I112").Value = 377
$ws.Range("J112").Value = 1513.7931
$ws.Range("K112").Value = 1131
$ws.Range("L112").Value = 4541.379300000001
$ws.Range("M112").Value = -23
$ws.Range("N112").Value = -6757.379300000001

$ws.Range("H130").Value = 55556944
$ws.Range("I130").Value = 125001130
$ws.Range("K130").Value = 375003390
$ws.Range("M130").Value = -374998370

$ws.Range("H131").Value = 3271.1042
$ws.Range("J131").Value = 3616.3489
$ws.Range("L131").Value = 10849.0467
$ws.Range("N131").Value = -20929.0467

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2890.8333
$ws.Range("I80").Value = 2780.8333
$ws.Range("J80").Value = 3000.8333
$ws.Range("K80").Value = 2780.8333
$ws.Range("L80").Value = 3000.8333
$ws.Range("M80").Value = -1782.8333
$ws.Range("N80").Value = -4996.8333

$ws.Range("H83").Value = 2890.8333
$ws.Range("I83").Value = 2780.8333
$ws.Range("J83").Value = 3000.8333
$ws.Range("K83").Value = 13904.1665
$ws.Range("L83").Value = 15004.1665
$ws.Range("M83").Value = -8912.166499999999
$ws.Range("N83").Value = -24988.1665

$ws.Range("H132").Value = 2519.6667
$ws.Range("I132").Value = 1773.3
$ws.Range("J132").Value = 5007.5557
$ws.Range("K132").Value = 5319.9
$ws.Range("L132").Value = 15022.6671
$ws.Range("M132").Value = -2789.9
$ws.Range("N132").Value = -20082.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2796.5
$ws.Range("I68").Value = 2080
$ws.Range("J68").Value = 3308.2856
$ws.Range("K68").Value = 2080
$ws.Range("L68").Value = 3308.2856
$ws.Range("M68").Value = -1331
$ws.Range("N68").Value = -4806.2856

$ws.Range("H71").Value = 2796.5
$ws.Range("I71").Value = 2080
$ws.Range("J71").Value = 3308.2856
$ws.Range("K71").Value = 10400
$ws.Range("L71").Value = 16541.428
$ws.Range("M71").Value = -6656
$ws.Range("N71").Value = -24029.428

$ws.Range("H82").Value = 2337.375
$ws.Range("I82").Value = 1339.8
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 1339.8
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -978.8
$ws.Range("N82").Value = -4722

$ws.Range("H85").Value = 2337.375
$ws.Range("I85").Value = 1339.8
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 1339.8
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -91.79999999999995
$ws.Range("N85").Value = -6496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4996
$ws.Range("J62").Value = 4996
$ws.Range("L62").Value = 4996
$ws.Range("N62").Value = -6244

$ws.Range("H65").Value = 4996
$ws.Range("J65").Value = 4996
$ws.Range("L65").Value = 24980
$ws.Range("N65").Value = -31220

$ws.Range("H81").Value = 1023.2632
$ws.Range("I81").Value = 764.2
$ws.Range("J81").Value = 1311.1111
$ws.Range("K81").Value = 1528.4
$ws.Range("L81").Value = 2622.2222
$ws.Range("M81").Value = -467.4000000000001
$ws.Range("N81").Value = -4744.2222

$ws.Range("H84").Value = 1023.2632
$ws.Range("I84").Value = 764.2
$ws.Range("J84").Value = 1311.1111
$ws.Range("K84").Value = 7642
$ws.Range("L84").Value = 13111.111
$ws.Range("M84").Value = -2338
$ws.Range("N84").Value = -23719.111

$ws.Range("H136").Value = 2727.3809
$ws.Range("I136").Value = 678.24
$ws.Range("J136").Value = 5740.8237
$ws.Range("K136").Value = 2034.72
$ws.Range("L136").Value = 17222.4711
$ws.Range("M136").Value = 515.28
$ws.Range("N136").Value = -22322.4711
